$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 8885.857
$ws.Cells.Item(38, 10).Value = 14466.667
$ws.Cells.Item(38, 12).Value = 43400.001
$ws.Cells.Item(38, 14).Value = -44144.001

$ws.Cells.Item(62, 8).Value = 6249.75
$ws.Cells.Item(62, 10).Value = 6999.6
$ws.Cells.Item(62, 12).Value = 6999.6
$ws.Cells.Item(62, 14).Value = -8247.6

$ws.Cells.Item(65, 8).Value = 6249.75
$ws.Cells.Item(65, 10).Value = 6999.6
$ws.Cells.Item(65, 12).Value = 34998
$ws.Cells.Item(65, 14).Value = -41238

$ws.Cells.Item(86, 8).Value = 125003140
$ws.Cells.Item(86, 9).Value = 125003140
$ws.Cells.Item(86, 11).Value = 125003140
$ws.Cells.Item(86, 13).Value = -125002017

$ws.Cells.Item(89, 8).Value = 125003140
$ws.Cells.Item(89, 9).Value = 125003140
$ws.Cells.Item(89, 11).Value = 625015700
$ws.Cells.Item(89, 13).Value = -625010084

$ws.Cells.Item(113, 8).Value = 3996.75
$ws.Cells.Item(113, 9).Value = 3993
$ws.Cells.Item(113, 10).Value = 3998
$ws.Cells.Item(113, 11).Value = 3993
$ws.Cells.Item(113, 12).Value = 3998
$ws.Cells.Item(113, 13).Value = -739
$ws.Cells.Item(113, 14).Value = -10506

$ws.Cells.Item(116, 8).Value = 1282754.5
$ws.Cells.Item(116, 9).Value = 2984337
$ws.Cells.Item(116, 10).Value = 6567.625
$ws.Cells.Item(116, 11).Value = 2984337
$ws.Cells.Item(116, 12).Value = 6567.625
$ws.Cells.Item(116, 13).Value = -2980895
$ws.Cells.Item(116, 14).Value = -13451.625

$ws.Cells.Item(132, 8).Value = 353858.72
$ws.Cells.Item(132, 9).Value = 597420.2
$ws.Cells.Item(132, 11).Value = 1792260.6
$ws.Cells.Item(132, 13).Value = -1789730.6

$ws.Cells.Item(138, 8).Value = 4308.1313
$ws.Cells.Item(138, 10).Value = 4869.9473
$ws.Cells.Item(138, 12).Value = 14609.8419
$ws.Cells.Item(138, 14).Value = -24889.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6717.436
$ws.Cells.Item(61, 9).Value = 6551.484
$ws.Cells.Item(61, 10).Value = 7360.5
$ws.Cells.Item(61, 11).Value = 6551.484
$ws.Cells.Item(61, 12).Value = 7360.5
$ws.Cells.Item(61, 13).Value = -6339.484
$ws.Cells.Item(61, 14).Value = -7784.5

$ws.Cells.Item(132, 8).Value = 18162.37
$ws.Cells.Item(132, 9).Value = 23631.072
$ws.Cells.Item(132, 11).Value = 70893.216
$ws.Cells.Item(132, 13).Value = -68363.216

$ws.Cells.Item(136, 8).Value = 6717.436
$ws.Cells.Item(136, 9).Value = 6551.484
$ws.Cells.Item(136, 10).Value = 7360.5
$ws.Cells.Item(136, 11).Value = 19654.452
$ws.Cells.Item(136, 12).Value = 22081.5
$ws.Cells.Item(136, 13).Value = -17104.452
$ws.Cells.Item(136, 14).Value = -27181.5

$ws.Cells.Item(138, 8).Value = 100000
$ws.Cells.Item(138, 10).Value = 100000
$ws.Cells.Item(138, 12).Value = 100000
$ws.Cells.Item(138, 14).Value = -110280

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 4027.516
$ws.Cells.Item(99, 10).Value = 5962.2
$ws.Cells.Item(99, 12).Value = 5962.2
$ws.Cells.Item(99, 14).Value = -8958.200000000001

$ws.Cells.Item(134, 8).Value = 3733.5
$ws.Cells.Item(134, 9).Value = 3667.2222
$ws.Cells.Item(134, 11).Value = 11001.6666
$ws.Cells.Item(134, 13).Value = -8466.6666

$ws.Cells.Item(140, 8).Value = 135153.55
$ws.Cells.Item(140, 10).Value = 135153.55
$ws.Cells.Item(140, 12).Value = 135153.55
$ws.Cells.Item(140, 14).Value = -145513.55

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1664.5
$ws.Cells.Item(16, 9).Value = 1329.375
$ws.Cells.Item(16, 10).Value = 2334.75
$ws.Cells.Item(16, 11).Value = 1329.375
$ws.Cells.Item(16, 12).Value = 2334.75
$ws.Cells.Item(16, 13).Value = -1042.375
$ws.Cells.Item(16, 14).Value = -2908.75

$ws.Cells.Item(22, 8).Value = 525.8182
$ws.Cells.Item(22, 9).Value = 383.57144
$ws.Cells.Item(22, 11).Value = 383.57144
$ws.Cells.Item(22, 13).Value = -33.57144

$ws.Cells.Item(31, 8).Value = 32261152
$ws.Cells.Item(31, 9).Value = 38463824
$ws.Cells.Item(31, 10).Value = 7259.8
$ws.Cells.Item(31, 11).Value = 38463824
$ws.Cells.Item(31, 12).Value = 7259.8
$ws.Cells.Item(31, 13).Value = -38463529
$ws.Cells.Item(31, 14).Value = -7849.8

$ws.Cells.Item(34, 8).Value = 32261152
$ws.Cells.Item(34, 9).Value = 38463824
$ws.Cells.Item(34, 10).Value = 7259.8
$ws.Cells.Item(34, 11).Value = 38463824
$ws.Cells.Item(34, 12).Value = 7259.8
$ws.Cells.Item(34, 13).Value = -38463622
$ws.Cells.Item(34, 14).Value = -7663.8

$ws.Cells.Item(58, 8).Value = 2470
$ws.Cells.Item(58, 10).Value = 3000
$ws.Cells.Item(58, 12).Value = 3000
$ws.Cells.Item(58, 14).Value = -3406

$ws.Cells.Item(99, 8).Value = 7064.95
$ws.Cells.Item(99, 9).Value = 4242.7144
$ws.Cells.Item(99, 10).Value = 8584.615
$ws.Cells.Item(99, 11).Value = 4242.7144
$ws.Cells.Item(99, 12).Value = 8584.615
$ws.Cells.Item(99, 13).Value = -2744.7144
$ws.Cells.Item(99, 14).Value = -11580.615

$ws.Cells.Item(113, 8).Value = 1664.5
$ws.Cells.Item(113, 9).Value = 1329.375
$ws.Cells.Item(113, 10).Value = 2334.75
$ws.Cells.Item(113, 11).Value = 1329.375
$ws.Cells.Item(113, 12).Value = 2334.75
$ws.Cells.Item(113, 13).Value = 840.625
$ws.Cells.Item(113, 14).Value = -6674.75

$ws.Cells.Item(122, 8).Value = 4219.0435
$ws.Cells.Item(122, 9).Value = 2871.3
$ws.Cells.Item(122, 10).Value = 5255.769
$ws.Cells.Item(122, 11).Value = 8613.900000000001
$ws.Cells.Item(122, 12).Value = 15767.307
$ws.Cells.Item(122, 13).Value = -6163.900000000001
$ws.Cells.Item(122, 14).Value = -20667.307

$ws.Cells.Item(126, 8).Value = 7064.95
$ws.Cells.Item(126, 9).Value = 4242.7144
$ws.Cells.Item(126, 10).Value = 8584.615
$ws.Cells.Item(126, 11).Value = 12728.1432
$ws.Cells.Item(126, 12).Value = 25753.845
$ws.Cells.Item(126, 13).Value = -10258.1432
$ws.Cells.Item(126, 14).Value = -30693.845

$ws.Cells.Item(132, 8).Value = 2618.8823
$ws.Cells.Item(132, 9).Value = 2364.2727
$ws.Cells.Item(132, 11).Value = 7092.8181
$ws.Cells.Item(132, 13).Value = -4562.8181

$ws.Cells.Item(134, 8).Value = 2593.6
$ws.Cells.Item(134, 9).Value = 2744.5
$ws.Cells.Item(134, 10).Value = 1990
$ws.Cells.Item(134, 11).Value = 8233.5
$ws.Cells.Item(134, 12).Value = 5970
$ws.Cells.Item(134, 13).Value = -5698.5
$ws.Cells.Item(134, 14).Value = -11040

$ws.Cells.Item(136, 8).Value = 2470
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 14).Value = -14100

$ws.Cells.Item(141, 8).Value = 333814
$ws.Cells.Item(141, 10).Value = 367572.4
$ws.Cells.Item(141, 12).Value = 367572.4
$ws.Cells.Item(141, 14).Value = -377932.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 3015.5
$ws.Cells.Item(122, 9).Value = 797.5
$ws.Cells.Item(122, 10).Value = 4124.5
$ws.Cells.Item(122, 11).Value = 7177.5
$ws.Cells.Item(122, 12).Value = 37120.5
$ws.Cells.Item(122, 13).Value = -4727.5
$ws.Cells.Item(122, 14).Value = -42020.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 34011460
$ws.Cells.Item(102, 9).Value = 56679876
$ws.Cells.Item(102, 11).Value = 56679876
$ws.Cells.Item(102, 13).Value = -56678254

$ws.Cells.Item(113, 8).Value = 1432.6666
$ws.Cells.Item(113, 9).Value = 1474.25
$ws.Cells.Item(113, 10).Value = 1399.4
$ws.Cells.Item(113, 11).Value = 1474.25
$ws.Cells.Item(113, 12).Value = 1399.4
$ws.Cells.Item(113, 13).Value = 695.75
$ws.Cells.Item(113, 14).Value = -5739.4

$ws.Cells.Item(123, 8).Value = 57107.11
$ws.Cells.Item(123, 10).Value = 57107.11
$ws.Cells.Item(123, 12).Value = 57107.11
$ws.Cells.Item(123, 14).Value = -62007.11

$ws.Cells.Item(132, 8).Value = 4558.0293
$ws.Cells.Item(132, 9).Value = 4429.643
$ws.Cells.Item(132, 11).Value = 13288.929
$ws.Cells.Item(132, 13).Value = -10758.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5648.7
$ws.Cells.Item(16, 9).Value = 5318.5415
$ws.Cells.Item(16, 10).Value = 6969.3335
$ws.Cells.Item(16, 11).Value = 5318.5415
$ws.Cells.Item(16, 12).Value = 6969.3335
$ws.Cells.Item(16, 13).Value = -5148.5415
$ws.Cells.Item(16, 14).Value = -7309.3335

$ws.Cells.Item(61, 8).Value = 3209.8635
$ws.Cells.Item(61, 9).Value = 3282.375
$ws.Cells.Item(61, 11).Value = 3282.375
$ws.Cells.Item(61, 13).Value = -3080.375

$ws.Cells.Item(113, 8).Value = 3209.8635
$ws.Cells.Item(113, 9).Value = 3282.375
$ws.Cells.Item(113, 11).Value = 3282.375
$ws.Cells.Item(113, 13).Value = -1112.375

$ws.Cells.Item(136, 8).Value = 3764.476
$ws.Cells.Item(136, 9).Value = 3256.7646
$ws.Cells.Item(136, 11).Value = 9770.293799999999
$ws.Cells.Item(136, 13).Value = -7220.293799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 12072.889
$ws.Cells.Item(62, 9).Value = 11465.3125
$ws.Cells.Item(62, 10).Value = 12956.637
$ws.Cells.Item(62, 11).Value = 11465.3125
$ws.Cells.Item(62, 12).Value = 12956.637
$ws.Cells.Item(62, 13).Value = -10841.3125
$ws.Cells.Item(62, 14).Value = -14204.637

$ws.Cells.Item(65, 8).Value = 12072.889
$ws.Cells.Item(65, 9).Value = 11465.3125
$ws.Cells.Item(65, 10).Value = 12956.637
$ws.Cells.Item(65, 11).Value = 57326.5625
$ws.Cells.Item(65, 12).Value = 64783.185
$ws.Cells.Item(65, 13).Value = -54206.5625
$ws.Cells.Item(65, 14).Value = -71023.185

$ws.Cells.Item(97, 8).Value = 99085
$ws.Cells.Item(97, 10).Value = 99085
$ws.Cells.Item(97, 12).Value = 99085
$ws.Cells.Item(97, 14).Value = -101067

$ws.Cells.Item(132, 8).Value = 6413698.5
$ws.Cells.Item(132, 9).Value = 12348199
$ws.Cells.Item(132, 11).Value = 37044597
$ws.Cells.Item(132, 13).Value = -37042067

$ws.Cells.Item(136, 8).Value = 5691.696
$ws.Cells.Item(136, 9).Value = 4200.7646
$ws.Cells.Item(136, 11).Value = 12602.2938
$ws.Cells.Item(136, 13).Value = -10052.2938
